# Paychart update:
#  - Eduardo Farias row: replace the lone "_GoBack" bookmark (left over at the
#    end of his dates cell) with a new "22/5" date run.
#  - Jorge Marin row: append a new ", 22/5" date run to his dates cell, and
#    move the "_GoBack" bookmark there (this is where Word now thinks the
#    last edit happened).
#
# NOTE: this runtime's Range/Selection.LanguageID setter does not persist to
# the underlying <w:rPr><w:lang .../></w:rPr> XML (verified experimentally -
# it silently no-ops, or in some code paths stamps the wrong paragraph
# entirely). Range.InsertXML, however, does round-trip <w:lang> correctly, so
# the new runs are authored as OOXML fragments and spliced in with
# InsertXML against a Range collapsed to the end of the target paragraph
# (i.e. right before its paragraph mark), which replaces that paragraph's
# content in place without touching sibling paragraphs/cells.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgClose = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Eduardo Farias: ", 15/5, " + bookmark  ->  ", 15/5, " + "22/5" run ---
$eduardoCell = $t.Cell(2, 2)
$eduardoXml = $pkgOpen + `
  '<w:p w:rsidR="008E4C32" w:rsidRPr="003B08F7" w:rsidRDefault="005E6A9B">' + `
    '<w:pPr>' + `
      '<w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>' + `
      '<w:rPr><w:lang w:val="en-US"/></w:rPr>' + `
    '</w:pPr>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>3/5</w:t></w:r>' + `
    '<w:r w:rsidR="00F97566"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, 8/5, 10/5</w:t></w:r>' + `
    '<w:r w:rsidR="003B08F7"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, 15/5, </w:t></w:r>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>22/5</w:t></w:r>' + `
  '</w:p>' + `
  $pkgClose

$wdCollapseEnd = [Microsoft.Office.Interop.Word.WdCollapseDirection]::wdCollapseEnd

$eduardoRange = $eduardoCell.Range
$eduardoRange.Collapse($wdCollapseEnd)
$eduardoRange.InsertXML($eduardoXml)

# --- Jorge Marin: ", 19/5"  ->  ", 19/5" + ", 22/5" run + moved _GoBack bookmark ---
$jorgeCell = $t.Cell(3, 2)
$jorgeXml = $pkgOpen + `
  '<w:p w:rsidR="007A0B59" w:rsidRDefault="00F97566">' + `
    '<w:pPr>' + `
      '<w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>' + `
      '<w:rPr><w:lang w:val="en-US"/></w:rPr>' + `
    '</w:pPr>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>8/5</w:t></w:r>' + `
    '<w:r w:rsidR="003B08F7"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>' + `
    '<w:r w:rsidR="003B08F7" w:rsidRPr="003B08F7"><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-US"/></w:rPr><w:t>12/5</w:t></w:r>' + `
    '<w:r w:rsidR="003B08F7"><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>' + `
    '<w:r w:rsidR="003B08F7" w:rsidRPr="003B08F7"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>15/5</w:t></w:r>' + `
    '<w:r w:rsidR="003B08F7"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, 19/5</w:t></w:r>' + `
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, 22/5</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
  '</w:p>' + `
  $pkgClose

$jorgeRange = $jorgeCell.Range
$jorgeRange.Collapse($wdCollapseEnd)
$jorgeRange.InsertXML($jorgeXml)

Write-Host "Eduardo cell now:" $eduardoCell.Range.Text
Write-Host "Jorge cell now:" $jorgeCell.Range.Text
